$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns: S1, T1, U1
$ws.Range("S1").Value = "Allow FS Deletion (Yes/No)"
$ws.Range("T1").Value = "Event Script"
$ws.Range("U1").Value = "Event Script Args"

# Match bold style of other header cells
$ws.Range("S1:U1").Font.Bold = $true

# Change A2 from "OCI Trial" to "First Flow"
$ws.Range("A2").Value = "First Flow"

# Column widths (best-fit recalculated by Excel after the content edits)
$ws.Columns("A").ColumnWidth = 16.109375
$ws.Columns("S").ColumnWidth = 23.6640625
$ws.Columns("T").ColumnWidth = 11
$ws.Columns("U").ColumnWidth = 15.21875
